$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06754433333333333
$ws.Range("H2").Value = 0.202633
$ws.Range("I2").Value = 0.006855017925354449
$ws.Range("J2").Value = 0.006855017925354449
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2679253333333333
$ws.Range("N2").Value = 0.8037759999999999
$ws.Range("O2").Value = 0.1226600350746756
$ws.Range("P2").Value = 0.1226600350746756
$ws.Range("Q2").Value = 0.01809683802311111
$ws.Range("R2").Value = 0.162871542208
$ws.Range("S2").Value = 0.0008408367391615066
$ws.Range("T2").Value = 0.0008408367391615068
# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06754433333333333
$ws.Range("H3").Value = 0.202633
$ws.Range("I3").Value = 0.006855017925354449
$ws.Range("J3").Value = 0.006855017925354449
$ws.Range("O3").Value = 0.327101565785771
$ws.Range("P3").Value = 0.327101565785771
$ws.Range("Q3").Value = 0.04825943551644444
$ws.Range("R3").Value = 0.434334919648
$ws.Range("S3").Value = 0.002242287096872967
$ws.Range("T3").Value = 0.002242287096872968
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06754433333333333
$ws.Range("H4").Value = 0.202633
$ws.Range("I4").Value = 0.006855017925354449
$ws.Range("J4").Value = 0.006855017925354449
$ws.Range("M4").Value = 1.145196333333333
$ws.Range("N4").Value = 3.435589
$ws.Range("O4").Value = 0.5242871984759059
$ws.Range("P4").Value = 0.5242871984759059
$ws.Range("Q4").Value = 0.07735152287077779
$ws.Range("R4").Value = 0.696163705837
$ws.Range("S4").Value = 0.0035939981435862
$ws.Range("T4").Value = 0.0035939981435862
# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.06754433333333333
$ws.Range("H5").Value = 0.202633
$ws.Range("I5").Value = 0.006855017925354449
$ws.Range("J5").Value = 0.006855017925354449
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.05668500000000001
$ws.Range("N5").Value = 0.170055
$ws.Range("O5").Value = 0.02595120066364754
$ws.Range("P5").Value = 0.02595120066364754
$ws.Range("Q5").Value = 0.003828750535
$ws.Range("R5").Value = 0.03445875481500001
$ws.Range("S5").Value = 0.0001778959457337741
$ws.Range("T5").Value = 0.0001778959457337741
# Row 6
$ws.Range("I6").Value = 0.7774992501642265
$ws.Range("J6").Value = 0.7774992501642265
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2679253333333333
$ws.Range("N6").Value = 0.8037759999999999
$ws.Range("O6").Value = 0.1226600350746756
$ws.Range("P6").Value = 0.1226600350746756
$ws.Range("Q6").Value = 2.052551597461334
$ws.Range("R6").Value = 18.472964377152
$ws.Range("S6").Value = 0.09536808529567799
$ws.Range("T6").Value = 0.09536808529567802
# Row 7
$ws.Range("I7").Value = 0.7774992501642265
$ws.Range("J7").Value = 0.7774992501642265
$ws.Range("O7").Value = 0.327101565785771
$ws.Range("P7").Value = 0.327101565785771
$ws.Range("S7").Value = 0.2543212221259813
$ws.Range("T7").Value = 0.2543212221259814
# Row 8
$ws.Range("I8").Value = 0.7774992501642265
$ws.Range("J8").Value = 0.7774992501642265
$ws.Range("M8").Value = 1.145196333333333
$ws.Range("N8").Value = 3.435589
$ws.Range("O8").Value = 0.5242871984759059
$ws.Range("P8").Value = 0.5242871984759059
$ws.Range("Q8").Value = 8.773244896800335
$ws.Range("R8").Value = 78.95920407120302
$ws.Range("S8").Value = 0.4076329036857198
$ws.Range("T8").Value = 0.4076329036857198
# Row 9
$ws.Range("I9").Value = 0.7774992501642265
$ws.Range("J9").Value = 0.7774992501642265
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.05668500000000001
$ws.Range("N9").Value = 0.170055
$ws.Range("O9").Value = 0.02595120066364754
$ws.Range("P9").Value = 0.02595120066364754
$ws.Range("Q9").Value = 0.4342586266650001
$ws.Range("R9").Value = 3.908327639985001
$ws.Range("S9").Value = 0.02017703905684734
$ws.Range("T9").Value = 0.02017703905684734
# Row 10
$ws.Range("G10").Value = 1.941983333333333
$ws.Range("H10").Value = 5.825949999999999
$ws.Range("I10").Value = 0.1970902650714284
$ws.Range("J10").Value = 0.1970902650714283
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2679253333333333
$ws.Range("N10").Value = 0.8037759999999999
$ws.Range("O10").Value = 0.1226600350746756
$ws.Range("P10").Value = 0.1226600350746756
$ws.Range("Q10").Value = 0.520306531911111
$ws.Range("R10").Value = 4.682758787199998
$ws.Range("S10").Value = 0.02417509882653851
$ws.Range("T10").Value = 0.02417509882653851
# Row 11
$ws.Range("G11").Value = 1.941983333333333
$ws.Range("H11").Value = 5.825949999999999
$ws.Range("I11").Value = 0.1970902650714284
$ws.Range("J11").Value = 0.1970902650714283
$ws.Range("O11").Value = 0.327101565785771
$ws.Range("P11").Value = 0.327101565785771
$ws.Range("Q11").Value = 1.387518609244444
$ws.Range("R11").Value = 12.4876674832
$ws.Range("S11").Value = 0.06446853430599687
$ws.Range("T11").Value = 0.06446853430599687
# Row 12
$ws.Range("G12").Value = 1.941983333333333
$ws.Range("H12").Value = 5.825949999999999
$ws.Range("I12").Value = 0.1970902650714284
$ws.Range("J12").Value = 0.1970902650714283
$ws.Range("M12").Value = 1.145196333333333
$ws.Range("N12").Value = 3.435589
$ws.Range("O12").Value = 0.5242871984759059
$ws.Range("P12").Value = 0.5242871984759059
$ws.Range("Q12").Value = 2.223952192727777
$ws.Range("R12").Value = 20.01556973455
$ws.Range("S12").Value = 0.1033319029211729
$ws.Range("T12").Value = 0.1033319029211728
# Row 13
$ws.Range("G13").Value = 1.941983333333333
$ws.Range("H13").Value = 5.825949999999999
$ws.Range("I13").Value = 0.1970902650714284
$ws.Range("J13").Value = 0.1970902650714283
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.05668500000000001
$ws.Range("N13").Value = 0.170055
$ws.Range("O13").Value = 0.02595120066364754
$ws.Range("P13").Value = 0.02595120066364754
$ws.Range("Q13").Value = 0.11008132525
$ws.Range("R13").Value = 0.9907319272499999
$ws.Range("S13").Value = 0.005114729017720121
$ws.Range("T13").Value = 0.005114729017720121
# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.182832
$ws.Range("H14").Value = 0.548496
$ws.Range("I14").Value = 0.01855546683899075
$ws.Range("J14").Value = 0.01855546683899075
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.2679253333333333
$ws.Range("N14").Value = 0.8037759999999999
$ws.Range("O14").Value = 0.1226600350746756
$ws.Range("P14").Value = 0.1226600350746756
$ws.Range("Q14").Value = 0.04898532454399999
$ws.Range("R14").Value = 0.440867920896
$ws.Range("S14").Value = 0.002276014213297586
$ws.Range("T14").Value = 0.002276014213297586
# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.182832
$ws.Range("H15").Value = 0.548496
$ws.Range("I15").Value = 0.01855546683899075
$ws.Range("J15").Value = 0.01855546683899075
$ws.Range("O15").Value = 0.327101565785771
$ws.Range("P15").Value = 0.327101565785771
$ws.Range("Q15").Value = 0.130630782464
$ws.Range("R15").Value = 1.175677042176
$ws.Range("S15").Value = 0.006069522256919826
$ws.Range("T15").Value = 0.006069522256919827
# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.182832
$ws.Range("H16").Value = 0.548496
$ws.Range("I16").Value = 0.01855546683899075
$ws.Range("J16").Value = 0.01855546683899075
$ws.Range("M16").Value = 1.145196333333333
$ws.Range("N16").Value = 3.435589
$ws.Range("O16").Value = 0.5242871984759059
$ws.Range("P16").Value = 0.5242871984759059
$ws.Range("Q16").Value = 0.209378536016
$ws.Range("R16").Value = 1.884406824144
$ws.Range("S16").Value = 0.009728393725427036
$ws.Range("T16").Value = 0.009728393725427036
# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.182832
$ws.Range("H17").Value = 0.548496
$ws.Range("I17").Value = 0.01855546683899075
$ws.Range("J17").Value = 0.01855546683899075
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.05668500000000001
$ws.Range("N17").Value = 0.170055
$ws.Range("O17").Value = 0.02595120066364754
$ws.Range("P17").Value = 0.02595120066364754
$ws.Range("Q17").Value = 0.01036383192
$ws.Range("R17").Value = 0.09327448728
$ws.Range("S17").Value = 0.0004815366433463067
$ws.Range("T17").Value = 0.0004815366433463067
